$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated coin/price/volume figures per the latest symbol-list refresh.
# Cells in columns D (Price) and E (Volume(1h)) hold numeric-looking text
# (e.g. "303.33", "-0.12%"), so force text formatting before assigning so
# Excel doesn't silently coerce them into numbers, then drop the temporary
# text format again so the cell style matches the rest of the sheet.
$updates = @(
    @{ Cell = 'D2'; Value = '303.33' },
    @{ Cell = 'E2'; Value = '-0.12%' },
    @{ Cell = 'D3'; Value = '37.05' },
    @{ Cell = 'E3'; Value = '2.48%' },
    @{ Cell = 'D4'; Value = '5.031' },
    @{ Cell = 'E4'; Value = '-0.70%' },
    @{ Cell = 'D5'; Value = '0.07831' },
    @{ Cell = 'D6'; Value = '2.209' },
    @{ Cell = 'E6'; Value = '-4.03%' },
    @{ Cell = 'D7'; Value = '7.991' },
    @{ Cell = 'E7'; Value = '-1.11%' },
    @{ Cell = 'D8'; Value = '0.9228' },
    @{ Cell = 'E8'; Value = '-0.33%' },
    @{ Cell = 'D9'; Value = '0.09879' },
    @{ Cell = 'E9'; Value = '-2.20%' },
    @{ Cell = 'D10'; Value = '0.1877' },
    @{ Cell = 'E10'; Value = '3.08%' },
    @{ Cell = 'D11'; Value = '0.08597' },
    @{ Cell = 'E11'; Value = '0.88%' },
    @{ Cell = 'D12'; Value = '0.03619' },
    @{ Cell = 'E12'; Value = '6.78%' },
    @{ Cell = 'D13'; Value = '0.09938' },
    @{ Cell = 'E13'; Value = '0.24%' },
    @{ Cell = 'D14'; Value = '0.001490' },
    @{ Cell = 'E14'; Value = '0.71%' },
    @{ Cell = 'D15'; Value = '0.005727' },
    @{ Cell = 'E15'; Value = '2.21%' },
    @{ Cell = 'D16'; Value = '3.459' },
    @{ Cell = 'D17'; Value = '4.025' },
    @{ Cell = 'E17'; Value = '0.72%' },
    @{ Cell = 'D18'; Value = '2.346' },
    @{ Cell = 'E18'; Value = '11.90%' },
    @{ Cell = 'D19'; Value = '0.3434' },
    @{ Cell = 'E19'; Value = '0.03%' },
    @{ Cell = 'D20'; Value = '0.1327' },
    @{ Cell = 'E20'; Value = '0.56%' },
    @{ Cell = 'D21'; Value = '4.773' },
    @{ Cell = 'E21'; Value = '4.75%' },
    @{ Cell = 'D22'; Value = '0.2202' },
    @{ Cell = 'E22'; Value = '-0.90%' },
    @{ Cell = 'D23'; Value = '0.04610' },
    @{ Cell = 'E23'; Value = '-1.31%' },
    @{ Cell = 'B24'; Value = 'HotbitToken' },
    @{ Cell = 'C24'; Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb' },
    @{ Cell = 'D24'; Value = '0.005196' },
    @{ Cell = 'E24'; Value = '15.78%' },
    @{ Cell = 'B25'; Value = 'BitKan' },
    @{ Cell = 'C25'; Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan' },
    @{ Cell = 'D25'; Value = '0.001251' },
    @{ Cell = 'E25'; Value = '0.87%' },
    @{ Cell = 'D26'; Value = '0.0001402' },
    @{ Cell = 'E26'; Value = '7.83%' },
    @{ Cell = 'D27'; Value = '0.0002719' },
    @{ Cell = 'E27'; Value = '-9.37%' },
    @{ Cell = 'D39'; Value = '0.01842' },
    @{ Cell = 'E39'; Value = '5.64%' },
    @{ Cell = 'E40'; Value = '1.49%' },
    @{ Cell = 'D41'; Value = '0.007954' },
    @{ Cell = 'E41'; Value = '1.80%' },
    @{ Cell = 'D42'; Value = '0.1403' },
    @{ Cell = 'E42'; Value = '-0.93%' },
    @{ Cell = 'D43'; Value = '0.007523' },
    @{ Cell = 'E43'; Value = '-14.55%' },
    @{ Cell = 'D44'; Value = '0.002243' },
    @{ Cell = 'E44'; Value = '12.17%' },
    @{ Cell = 'D45'; Value = '0.01042' },
    @{ Cell = 'E45'; Value = '13.51%' },
    @{ Cell = 'D46'; Value = '0.00006293' },
    @{ Cell = 'E46'; Value = '4.74%' },
    @{ Cell = 'D47'; Value = '0.00000000751' },
    @{ Cell = 'E47'; Value = '0.14%' },
    @{ Cell = 'D48'; Value = '0.0005803' },
    @{ Cell = 'E48'; Value = '0.04%' },
    @{ Cell = 'D49'; Value = '35.80' },
    @{ Cell = 'E49'; Value = '516.85%' },
    @{ Cell = 'D50'; Value = '0.002690' },
    @{ Cell = 'E50'; Value = '0.03%' },
    @{ Cell = 'D51'; Value = '0.00002103' },
    @{ Cell = 'E51'; Value = '0.14%' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $col = $u.Cell -replace '[0-9]+$', ''
    if ($col -eq 'D' -or $col -eq 'E') {
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        $cell.ClearFormats()
    } else {
        $cell.Value = $u.Value
    }
}
